# "Perfil jugador, historial y config"
#
# Hoja3 (database-schema notes sheet) is reworked:
#  - the "match_hist" table header becomes "match_det"
#  - match_status domain note is expanded from "jugado|pendiente" to
#    "pendiente | confirmado | jugado"
#  - the old user_cal / FLOAT(2,2) / calification rating fields are dropped
#    and replaced by user_pydmchs / user_wonmatches / user_lostmatches
#    (all INT(10) "played/won/lost matches" counters)
#  - a little 3-row "workflow" legend is appended below the existing notes
#  - a few column widths are widened to fit the new text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- header rename: match_hist -> match_det -------------------------------
$ws.Range("H1").Value = "match_det"

# --- match_status note: add the "confirmado" state ------------------------
$ws.Range("F5").Value = "pendiente | confirmado | jugado"

# --- new workflow legend: rows 20-22 ---------------------------------------
$ws.Range("D20").Value = "crear ->"
$ws.Range("E20").Value = "pendiente"

$ws.Range("D21").Value = "10 jugadores ->"
$ws.Range("E21").Value = "confirmado"

$ws.Range("D22").Value = "cargar result ->"
$ws.Range("E22").Value = "jugado"

# --- replace the calification trio with a won/lost-matches trio -----------
$ws.Range("A9").Value = "INT(10)"
$ws.Range("B9").Value = "user_pydmchs"
$ws.Range("C9").Value = "played matches"

$ws.Range("B10").Value = "user_wonmatches"
$ws.Range("C10").ClearContents()

$ws.Range("B11").Value = "user_lostmatches"

# --- widen a few columns to fit the new text -------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.28515625
$ws.Columns.Item(4).ColumnWidth = 15.85546875
$ws.Columns.Item(6).ColumnWidth = 33.42578125

# --- move the view/selection down to the new rows --------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
